# Fruta / hortaliza, semanal
# Insert a new weekly observation as row 11, shifting existing data down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11 (pushes rows 11..40 down to 12..41,
# copying formatting from the row above, same as Excel's default behaviour).
$ws.Rows("11:11").Insert()

# Populate the newly inserted row 11 with the new weekly record.
$ws.Cells.Item(11, 1).Value2 = 1
$ws.Cells.Item(11, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(11, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(11, 4).Value2 = 45133
$ws.Cells.Item(11, 5).Value2 = 15
$ws.Cells.Item(11, 6).Value2 = 100112013
$ws.Cells.Item(11, 7).Value2 = "Alcachofa"
$ws.Cells.Item(11, 8).Value2 = "Madrigal"
$ws.Cells.Item(11, 9).Value2 = "Primera"
$ws.Cells.Item(11, 10).Value2 = 170
$ws.Cells.Item(11, 11).Value2 = 20000
$ws.Cells.Item(11, 12).Value2 = 21000
$ws.Cells.Item(11, 13).Value2 = 20529
$ws.Cells.Item(11, 14).Value2 = "$/caja 40 unidades"
$ws.Cells.Item(11, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(11, 16).Value2 = 513
$ws.Cells.Item(11, 17).Value2 = 40
$ws.Cells.Item(11, 18).Value2 = "Hortaliza"
